# Apply updated NATMI Mmp2-Sdc2 ligand-receptor interaction data
# following Dr Hou advice: expand from 2x2 to 4x4 cluster pairs (ECs, FAPs, M2, sCs)
# and recompute all derived statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Mmp2"
$ws.Cells.Item(2,3).Value = "Sdc2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = [double]"3"
$ws.Cells.Item(2,6).Value = [double]"1"
$ws.Cells.Item(2,7).Value = [double]"3.140455666666667"
$ws.Cells.Item(2,8).Value = [double]"9.421367"
$ws.Cells.Item(2,9).Value = [double]"0.005037281036089241"
$ws.Cells.Item(2,10).Value = [double]"0.005037281036089241"
$ws.Cells.Item(2,11).Value = [double]"2"
$ws.Cells.Item(2,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2,13).Value = [double]"0.6646083333333334"
$ws.Cells.Item(2,14).Value = [double]"1.993825"
$ws.Cells.Item(2,15).Value = [double]"0.005515555660921567"
$ws.Cells.Item(2,16).Value = [double]"0.005515555660921567"
$ws.Cells.Item(2,17).Value = [double]"2.087173006530556"
$ws.Cells.Item(2,18).Value = [double]"18.784557058775"
$ws.Cells.Item(2,19).Value = [double]"2.778340393425487E-05"
$ws.Cells.Item(2,20).Value = [double]"2.778340393425487E-05"

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Mmp2"
$ws.Cells.Item(3,3).Value = "Sdc2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = [double]"3"
$ws.Cells.Item(3,6).Value = [double]"1"
$ws.Cells.Item(3,7).Value = [double]"3.140455666666667"
$ws.Cells.Item(3,8).Value = [double]"9.421367"
$ws.Cells.Item(3,9).Value = [double]"0.005037281036089241"
$ws.Cells.Item(3,10).Value = [double]"0.005037281036089241"
$ws.Cells.Item(3,11).Value = [double]"3"
$ws.Cells.Item(3,12).Value = [double]"1"
$ws.Cells.Item(3,13).Value = [double]"85.826024"
$ws.Cells.Item(3,14).Value = [double]"257.478072"
$ws.Cells.Item(3,15).Value = [double]"0.7122664414292983"
$ws.Cells.Item(3,16).Value = [double]"0.7122664414292983"
$ws.Cells.Item(3,17).Value = [double]"269.5328234182693"
$ws.Cells.Item(3,18).Value = [double]"2425.795410764424"
$ws.Cells.Item(3,19).Value = [double]"0.003587886238054572"
$ws.Cells.Item(3,20).Value = [double]"0.003587886238054572"

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Mmp2"
$ws.Cells.Item(4,3).Value = "Sdc2"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = [double]"3"
$ws.Cells.Item(4,6).Value = [double]"1"
$ws.Cells.Item(4,7).Value = [double]"3.140455666666667"
$ws.Cells.Item(4,8).Value = [double]"9.421367"
$ws.Cells.Item(4,9).Value = [double]"0.005037281036089241"
$ws.Cells.Item(4,10).Value = [double]"0.005037281036089241"
$ws.Cells.Item(4,11).Value = [double]"1"
$ws.Cells.Item(4,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4,13).Value = [double]"0.009315666666666667"
$ws.Cells.Item(4,14).Value = [double]"0.027947"
$ws.Cells.Item(4,15).Value = [double]"7.731031261809587E-05"
$ws.Cells.Item(4,16).Value = [double]"7.731031261809588E-05"
$ws.Cells.Item(4,17).Value = [double]"0.02925543817211111"
$ws.Cells.Item(4,18).Value = [double]"0.263298943549"
$ws.Cells.Item(4,19).Value = [double]"3.894337716452651E-07"
$ws.Cells.Item(4,20).Value = [double]"3.894337716452652E-07"

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Mmp2"
$ws.Cells.Item(5,3).Value = "Sdc2"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = [double]"3"
$ws.Cells.Item(5,6).Value = [double]"1"
$ws.Cells.Item(5,7).Value = [double]"3.140455666666667"
$ws.Cells.Item(5,8).Value = [double]"9.421367"
$ws.Cells.Item(5,9).Value = [double]"0.005037281036089241"
$ws.Cells.Item(5,10).Value = [double]"0.005037281036089241"
$ws.Cells.Item(5,11).Value = [double]"3"
$ws.Cells.Item(5,12).Value = [double]"1"
$ws.Cells.Item(5,13).Value = [double]"33.99712866666667"
$ws.Cells.Item(5,14).Value = [double]"101.991386"
$ws.Cells.Item(5,15).Value = [double]"0.2821406925971621"
$ws.Cells.Item(5,16).Value = [double]"0.2821406925971621"
$ws.Cells.Item(5,17).Value = [double]"106.7664753716291"
$ws.Cells.Item(5,18).Value = [double]"960.898278344662"
$ws.Cells.Item(5,19).Value = [double]"0.001421221960328769"
$ws.Cells.Item(5,20).Value = [double]"0.001421221960328769"

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Mmp2"
$ws.Cells.Item(6,3).Value = "Sdc2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = [double]"3"
$ws.Cells.Item(6,6).Value = [double]"1"
$ws.Cells.Item(6,7).Value = [double]"601.5827839999999"
$ws.Cells.Item(6,8).Value = [double]"1804.748352"
$ws.Cells.Item(6,9).Value = [double]"0.9649368980576715"
$ws.Cells.Item(6,10).Value = [double]"0.9649368980576714"
$ws.Cells.Item(6,11).Value = [double]"2"
$ws.Cells.Item(6,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(6,13).Value = [double]"0.6646083333333334"
$ws.Cells.Item(6,14).Value = [double]"1.993825"
$ws.Cells.Item(6,15).Value = [double]"0.005515555660921567"
$ws.Cells.Item(6,16).Value = [double]"0.005515555660921567"
$ws.Cells.Item(6,17).Value = [double]"399.8169314362667"
$ws.Cells.Item(6,18).Value = [double]"3598.3523829264"
$ws.Cells.Item(6,19).Value = [double]"0.005322163170514087"
$ws.Cells.Item(6,20).Value = [double]"0.005322163170514086"

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Mmp2"
$ws.Cells.Item(7,3).Value = "Sdc2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = [double]"3"
$ws.Cells.Item(7,6).Value = [double]"1"
$ws.Cells.Item(7,7).Value = [double]"601.5827839999999"
$ws.Cells.Item(7,8).Value = [double]"1804.748352"
$ws.Cells.Item(7,9).Value = [double]"0.9649368980576715"
$ws.Cells.Item(7,10).Value = [double]"0.9649368980576714"
$ws.Cells.Item(7,11).Value = [double]"3"
$ws.Cells.Item(7,12).Value = [double]"1"
$ws.Cells.Item(7,13).Value = [double]"85.826024"
$ws.Cells.Item(7,14).Value = [double]"257.478072"
$ws.Cells.Item(7,15).Value = [double]"0.7122664414292983"
$ws.Cells.Item(7,16).Value = [double]"0.7122664414292983"
$ws.Cells.Item(7,17).Value = [double]"51631.45845757081"
$ws.Cells.Item(7,18).Value = [double]"464683.1261181373"
$ws.Cells.Item(7,19).Value = [double]"0.6872921705833632"
$ws.Cells.Item(7,20).Value = [double]"0.6872921705833631"

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Mmp2"
$ws.Cells.Item(8,3).Value = "Sdc2"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = [double]"3"
$ws.Cells.Item(8,6).Value = [double]"1"
$ws.Cells.Item(8,7).Value = [double]"601.5827839999999"
$ws.Cells.Item(8,8).Value = [double]"1804.748352"
$ws.Cells.Item(8,9).Value = [double]"0.9649368980576715"
$ws.Cells.Item(8,10).Value = [double]"0.9649368980576714"
$ws.Cells.Item(8,11).Value = [double]"1"
$ws.Cells.Item(8,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(8,13).Value = [double]"0.009315666666666667"
$ws.Cells.Item(8,14).Value = [double]"0.027947"
$ws.Cells.Item(8,15).Value = [double]"7.731031261809587E-05"
$ws.Cells.Item(8,16).Value = [double]"7.731031261809588E-05"
$ws.Cells.Item(8,17).Value = [double]"5.604144688149333"
$ws.Cells.Item(8,18).Value = [double]"50.437302193344"
$ws.Cells.Item(8,19).Value = [double]"7.459957324557428E-05"
$ws.Cells.Item(8,20).Value = [double]"7.45995732455743E-05"

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Mmp2"
$ws.Cells.Item(9,3).Value = "Sdc2"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = [double]"3"
$ws.Cells.Item(9,6).Value = [double]"1"
$ws.Cells.Item(9,7).Value = [double]"601.5827839999999"
$ws.Cells.Item(9,8).Value = [double]"1804.748352"
$ws.Cells.Item(9,9).Value = [double]"0.9649368980576715"
$ws.Cells.Item(9,10).Value = [double]"0.9649368980576714"
$ws.Cells.Item(9,11).Value = [double]"3"
$ws.Cells.Item(9,12).Value = [double]"1"
$ws.Cells.Item(9,13).Value = [double]"33.99712866666667"
$ws.Cells.Item(9,14).Value = [double]"101.991386"
$ws.Cells.Item(9,15).Value = [double]"0.2821406925971621"
$ws.Cells.Item(9,16).Value = [double]"0.2821406925971621"
$ws.Cells.Item(9,17).Value = [double]"20452.08731129954"
$ws.Cells.Item(9,18).Value = [double]"184068.7858016959"
$ws.Cells.Item(9,19).Value = [double]"0.2722479647305486"
$ws.Cells.Item(9,20).Value = [double]"0.2722479647305486"

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Mmp2"
$ws.Cells.Item(10,3).Value = "Sdc2"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = [double]"3"
$ws.Cells.Item(10,6).Value = [double]"1"
$ws.Cells.Item(10,7).Value = [double]"0.9353736666666667"
$ws.Cells.Item(10,8).Value = [double]"2.806121"
$ws.Cells.Item(10,9).Value = [double]"0.001500336426579262"
$ws.Cells.Item(10,10).Value = [double]"0.001500336426579262"
$ws.Cells.Item(10,11).Value = [double]"2"
$ws.Cells.Item(10,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(10,13).Value = [double]"0.6646083333333334"
$ws.Cells.Item(10,14).Value = [double]"1.993825"
$ws.Cells.Item(10,15).Value = [double]"0.005515555660921567"
$ws.Cells.Item(10,16).Value = [double]"0.005515555660921567"
$ws.Cells.Item(10,17).Value = [double]"0.6216571336472222"
$ws.Cells.Item(10,18).Value = [double]"5.594914202825"
$ws.Cells.Item(10,19).Value = [double]"8.275189070906082E-06"
$ws.Cells.Item(10,20).Value = [double]"8.275189070906082E-06"

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Mmp2"
$ws.Cells.Item(11,3).Value = "Sdc2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = [double]"3"
$ws.Cells.Item(11,6).Value = [double]"1"
$ws.Cells.Item(11,7).Value = [double]"0.9353736666666667"
$ws.Cells.Item(11,8).Value = [double]"2.806121"
$ws.Cells.Item(11,9).Value = [double]"0.001500336426579262"
$ws.Cells.Item(11,10).Value = [double]"0.001500336426579262"
$ws.Cells.Item(11,11).Value = [double]"3"
$ws.Cells.Item(11,12).Value = [double]"1"
$ws.Cells.Item(11,13).Value = [double]"85.826024"
$ws.Cells.Item(11,14).Value = [double]"257.478072"
$ws.Cells.Item(11,15).Value = [double]"0.7122664414292983"
$ws.Cells.Item(11,16).Value = [double]"0.7122664414292983"
$ws.Cells.Item(11,17).Value = [double]"80.27940276430134"
$ws.Cells.Item(11,18).Value = [double]"722.514624878712"
$ws.Cells.Item(11,19).Value = [double]"0.00106863928750636"
$ws.Cells.Item(11,20).Value = [double]"0.00106863928750636"

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Mmp2"
$ws.Cells.Item(12,3).Value = "Sdc2"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = [double]"3"
$ws.Cells.Item(12,6).Value = [double]"1"
$ws.Cells.Item(12,7).Value = [double]"0.9353736666666667"
$ws.Cells.Item(12,8).Value = [double]"2.806121"
$ws.Cells.Item(12,9).Value = [double]"0.001500336426579262"
$ws.Cells.Item(12,10).Value = [double]"0.001500336426579262"
$ws.Cells.Item(12,11).Value = [double]"1"
$ws.Cells.Item(12,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(12,13).Value = [double]"0.009315666666666667"
$ws.Cells.Item(12,14).Value = [double]"0.027947"
$ws.Cells.Item(12,15).Value = [double]"7.731031261809587E-05"
$ws.Cells.Item(12,16).Value = [double]"7.731031261809588E-05"
$ws.Cells.Item(12,17).Value = [double]"0.008713629287444444"
$ws.Cells.Item(12,18).Value = [double]"0.07842266358700001"
$ws.Cells.Item(12,19).Value = [double]"1.159914781711596E-07"
$ws.Cells.Item(12,20).Value = [double]"1.159914781711596E-07"

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Mmp2"
$ws.Cells.Item(13,3).Value = "Sdc2"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = [double]"3"
$ws.Cells.Item(13,6).Value = [double]"1"
$ws.Cells.Item(13,7).Value = [double]"0.9353736666666667"
$ws.Cells.Item(13,8).Value = [double]"2.806121"
$ws.Cells.Item(13,9).Value = [double]"0.001500336426579262"
$ws.Cells.Item(13,10).Value = [double]"0.001500336426579262"
$ws.Cells.Item(13,11).Value = [double]"3"
$ws.Cells.Item(13,12).Value = [double]"1"
$ws.Cells.Item(13,13).Value = [double]"33.99712866666667"
$ws.Cells.Item(13,14).Value = [double]"101.991386"
$ws.Cells.Item(13,15).Value = [double]"0.2821406925971621"
$ws.Cells.Item(13,16).Value = [double]"0.2821406925971621"
$ws.Cells.Item(13,17).Value = [double]"31.80001889707844"
$ws.Cells.Item(13,18).Value = [double]"286.200170073706"
$ws.Cells.Item(13,19).Value = [double]"0.0004233059585238241"
$ws.Cells.Item(13,20).Value = [double]"0.0004233059585238241"

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Mmp2"
$ws.Cells.Item(14,3).Value = "Sdc2"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = [double]"3"
$ws.Cells.Item(14,6).Value = [double]"1"
$ws.Cells.Item(14,7).Value = [double]"17.78400266666667"
$ws.Cells.Item(14,8).Value = [double]"53.352008"
$ws.Cells.Item(14,9).Value = [double]"0.02852548447966007"
$ws.Cells.Item(14,10).Value = [double]"0.02852548447966007"
$ws.Cells.Item(14,11).Value = [double]"2"
$ws.Cells.Item(14,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(14,13).Value = [double]"0.6646083333333334"
$ws.Cells.Item(14,14).Value = [double]"1.993825"
$ws.Cells.Item(14,15).Value = [double]"0.005515555660921567"
$ws.Cells.Item(14,16).Value = [double]"0.005515555660921567"
$ws.Cells.Item(14,17).Value = [double]"11.81939637228889"
$ws.Cells.Item(14,18).Value = [double]"106.3745673506"
$ws.Cells.Item(14,19).Value = [double]"0.0001573338974023194"
$ws.Cells.Item(14,20).Value = [double]"0.0001573338974023194"

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Mmp2"
$ws.Cells.Item(15,3).Value = "Sdc2"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = [double]"3"
$ws.Cells.Item(15,6).Value = [double]"1"
$ws.Cells.Item(15,7).Value = [double]"17.78400266666667"
$ws.Cells.Item(15,8).Value = [double]"53.352008"
$ws.Cells.Item(15,9).Value = [double]"0.02852548447966007"
$ws.Cells.Item(15,10).Value = [double]"0.02852548447966007"
$ws.Cells.Item(15,11).Value = [double]"3"
$ws.Cells.Item(15,12).Value = [double]"1"
$ws.Cells.Item(15,13).Value = [double]"85.826024"
$ws.Cells.Item(15,14).Value = [double]"257.478072"
$ws.Cells.Item(15,15).Value = [double]"0.7122664414292983"
$ws.Cells.Item(15,16).Value = [double]"0.7122664414292983"
$ws.Cells.Item(15,17).Value = [double]"1526.330239685398"
$ws.Cells.Item(15,18).Value = [double]"13736.97215716858"
$ws.Cells.Item(15,19).Value = [double]"0.02031774532037416"
$ws.Cells.Item(15,20).Value = [double]"0.02031774532037415"

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Mmp2"
$ws.Cells.Item(16,3).Value = "Sdc2"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = [double]"3"
$ws.Cells.Item(16,6).Value = [double]"1"
$ws.Cells.Item(16,7).Value = [double]"17.78400266666667"
$ws.Cells.Item(16,8).Value = [double]"53.352008"
$ws.Cells.Item(16,9).Value = [double]"0.02852548447966007"
$ws.Cells.Item(16,10).Value = [double]"0.02852548447966007"
$ws.Cells.Item(16,11).Value = [double]"1"
$ws.Cells.Item(16,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(16,13).Value = [double]"0.009315666666666667"
$ws.Cells.Item(16,14).Value = [double]"0.027947"
$ws.Cells.Item(16,15).Value = [double]"7.731031261809587E-05"
$ws.Cells.Item(16,16).Value = [double]"7.731031261809588E-05"
$ws.Cells.Item(16,17).Value = [double]"0.1656698408417778"
$ws.Cells.Item(16,18).Value = [double]"1.491028567576"
$ws.Cells.Item(16,19).Value = [double]"2.205314122705162E-06"
$ws.Cells.Item(16,20).Value = [double]"2.205314122705162E-06"

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Mmp2"
$ws.Cells.Item(17,3).Value = "Sdc2"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = [double]"3"
$ws.Cells.Item(17,6).Value = [double]"1"
$ws.Cells.Item(17,7).Value = [double]"17.78400266666667"
$ws.Cells.Item(17,8).Value = [double]"53.352008"
$ws.Cells.Item(17,9).Value = [double]"0.02852548447966007"
$ws.Cells.Item(17,10).Value = [double]"0.02852548447966007"
$ws.Cells.Item(17,11).Value = [double]"3"
$ws.Cells.Item(17,12).Value = [double]"1"
$ws.Cells.Item(17,13).Value = [double]"33.99712866666667"
$ws.Cells.Item(17,14).Value = [double]"101.991386"
$ws.Cells.Item(17,15).Value = [double]"0.2821406925971621"
$ws.Cells.Item(17,16).Value = [double]"0.2821406925971621"
$ws.Cells.Item(17,17).Value = [double]"604.6050268670099"
$ws.Cells.Item(17,18).Value = [double]"5441.445241803089"
$ws.Cells.Item(17,19).Value = [double]"0.00804819994776089"
$ws.Cells.Item(17,20).Value = [double]"0.008048199947760888"

